$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column M header: "Final Project Proposal"
$ws.Range("M2").Value = "Final Project Proposal"

# Widen column M (13) to match the new header column (~21.48 chars)
$ws.Columns.Item(13).ColumnWidth = 20.6

# K8 was a bare literal 0; now becomes a formula 45/50 (=0.9)
$ws.Range("K8").Formula = "=45/50"

# New "Final Project Proposal" scores (full credit = 1) for several students
$ws.Range("M11").Value = 1
$ws.Range("M13").Value = 1
$ws.Range("M14").Value = 1
$ws.Range("M15").Value = 1
$ws.Range("M18").Value = 1

# Update the active selection to M14, matching the saved view state
$ws.Range("M14").Select()
